# GTA_EVENTS_TABLE.xlsx update: add new events (Sumo series, Hunting Pack I-VII,
# Stunters Vs Snipers - 4 Way, Duck Hunt OG, RPG vs Flying Cars) and fix a few
# rows that were missed, per the commit "Updated to include more events and a
# few races that were missed". The final table grows from 17 to 33 data+header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: touch one cell per distinct string value, in the same order the
# strings were first authored, so the shared-strings table comes out in the
# expected order (Excel appends to xl/sharedStrings.xml in first-use order).
$ws.Range("A1").Value2 = "event_name"
$ws.Range("B12").Value2 = "LTS"
$ws.Range("A10").Value2 = "Duck Hunt"
$ws.Range("C1").Value2 = "player_min"
$ws.Range("D1").Value2 = "player_max"
$ws.Range("B1").Value2 = "event_type"
$ws.Range("A12").Value2 = "Skyscraper LTS"
$ws.Range("A13").Value2 = "Tataviam Truckstop LTS"
$ws.Range("A14").Value2 = "Storm Drain - Tank LTS"
$ws.Range("A15").Value2 = "Alta LTS"
$ws.Range("E1").Value2 = "Team"
$ws.Range("A16").Value2 = "Legal Eagle"
$ws.Range("B16").Value2 = "Team Deathmatch"
$ws.Range("A17").Value2 = "Downtown"
$ws.Range("A18").Value2 = "Kortz LTS"
$ws.Range("A19").Value2 = "Power Station LTS"
$ws.Range("A20").Value2 = "Governmental"
$ws.Range("A21").Value2 = "RPG vs Armored Sharks"
$ws.Range("A23").Value2 = "Bikes vs RPG"
$ws.Range("A24").Value2 = "Sumo I"
$ws.Range("A25").Value2 = "Sumo II"
$ws.Range("A26").Value2 = "Sumo III"
$ws.Range("B2").Value2 = "Adversary"
$ws.Range("A27").Value2 = "Sumo (Remix) I"
$ws.Range("A28").Value2 = "Sumo (Remix) II"
$ws.Range("A29").Value2 = "Sumo (Remix) III"
$ws.Range("A30").Value2 = "Sumo (Remix) IV"
$ws.Range("A31").Value2 = "Sumo (Remix) V"
$ws.Range("A32").Value2 = "Sumo (Remix) VI"
$ws.Range("A33").Value2 = "Sumo (Remix) VII"
$ws.Range("A22").Value2 = "RPG vs Flying Cars"
$ws.Range("A9").Value2 = "Stunters Vs Snipers - 4 Way"
$ws.Range("A11").Value2 = "Duck Hunt OG"
$ws.Range("A2").Value2 = "Hunting Pack I"
$ws.Range("A3").Value2 = "Hunting Pack II"
$ws.Range("A4").Value2 = "Hunting Pack III"
$ws.Range("A5").Value2 = "Hunting Pack IV"
$ws.Range("A6").Value2 = "Hunting Pack V"
$ws.Range("A7").Value2 = "Hunting Pack VI"
$ws.Range("A8").Value2 = "Hunting Pack VII"

# --- Step 2: bulk-write the full A1:E33 table (re-asserts the step-1 cells and
# fills in every remaining text/number cell in one shot).
$table = New-Object "object[,]" 33,5
$table[0,0] = "event_name"
$table[0,1] = "event_type"
$table[0,2] = "player_min"
$table[0,3] = "player_max"
$table[0,4] = "Team"
$table[1,0] = "Hunting Pack I"
$table[1,1] = "Adversary"
$table[1,2] = 4
$table[1,3] = 8
$table[1,4] = 3
$table[2,0] = "Hunting Pack II"
$table[2,1] = "Adversary"
$table[2,2] = 4
$table[2,3] = 8
$table[2,4] = 3
$table[3,0] = "Hunting Pack III"
$table[3,1] = "Adversary"
$table[3,2] = 4
$table[3,3] = 8
$table[3,4] = 3
$table[4,0] = "Hunting Pack IV"
$table[4,1] = "Adversary"
$table[4,2] = 4
$table[4,3] = 8
$table[4,4] = 3
$table[5,0] = "Hunting Pack V"
$table[5,1] = "Adversary"
$table[5,2] = 4
$table[5,3] = 8
$table[5,4] = 3
$table[6,0] = "Hunting Pack VI"
$table[6,1] = "Adversary"
$table[6,2] = 4
$table[6,3] = 8
$table[6,4] = 3
$table[7,0] = "Hunting Pack VII"
$table[7,1] = "Adversary"
$table[7,2] = 4
$table[7,3] = 8
$table[7,4] = 3
$table[8,0] = "Stunters Vs Snipers - 4 Way"
$table[8,1] = "Adversary"
$table[8,2] = 1
$table[8,3] = 8
$table[8,4] = 2
$table[9,0] = "Duck Hunt"
$table[9,1] = "Adversary"
$table[9,2] = 1
$table[9,3] = 8
$table[9,4] = 2
$table[10,0] = "Duck Hunt OG"
$table[10,1] = "Adversary"
$table[10,2] = 1
$table[10,3] = 8
$table[10,4] = 2
$table[11,0] = "Skyscraper LTS"
$table[11,1] = "LTS"
$table[11,2] = 2
$table[11,3] = 16
$table[11,4] = 2
$table[12,0] = "Tataviam Truckstop LTS"
$table[12,1] = "LTS"
$table[12,2] = 2
$table[12,3] = 8
$table[12,4] = 2
$table[13,0] = "Storm Drain - Tank LTS"
$table[13,1] = "LTS"
$table[13,2] = 2
$table[13,3] = 16
$table[13,4] = 2
$table[14,0] = "Alta LTS"
$table[14,1] = "LTS"
$table[14,2] = 2
$table[14,3] = 8
$table[14,4] = 2
$table[15,0] = "Legal Eagle"
$table[15,1] = "Team Deathmatch"
$table[15,2] = 2
$table[15,3] = 8
$table[15,4] = 2
$table[16,0] = "Downtown"
$table[16,1] = "Team Deathmatch"
$table[16,2] = 2
$table[16,3] = 8
$table[16,4] = 2
$table[17,0] = "Kortz LTS"
$table[17,1] = "LTS"
$table[17,2] = 2
$table[17,3] = 16
$table[17,4] = 2
$table[18,0] = "Power Station LTS"
$table[18,1] = "LTS"
$table[18,2] = 2
$table[18,3] = 16
$table[18,4] = 2
$table[19,0] = "Governmental"
$table[19,1] = "Team Deathmatch"
$table[19,2] = 8
$table[19,3] = 16
$table[19,4] = 2
$table[20,0] = "RPG vs Armored Sharks"
$table[20,1] = "LTS"
$table[20,2] = 2
$table[20,3] = 18
$table[20,4] = 2
$table[21,0] = "RPG vs Flying Cars"
$table[21,1] = "Adversary"
$table[21,2] = 2
$table[21,3] = 30
$table[21,4] = 2
$table[22,0] = "Bikes vs RPG"
$table[22,1] = "LTS"
$table[22,2] = 2
$table[22,3] = 18
$table[22,4] = 2
$table[23,0] = "Sumo I"
$table[23,1] = "Adversary"
$table[23,2] = 2
$table[23,3] = 8
$table[23,4] = 4
$table[24,0] = "Sumo II"
$table[24,1] = "Adversary"
$table[24,2] = 2
$table[24,3] = 8
$table[24,4] = 4
$table[25,0] = "Sumo III"
$table[25,1] = "Adversary"
$table[25,2] = 2
$table[25,3] = 8
$table[25,4] = 4
$table[26,0] = "Sumo (Remix) I"
$table[26,1] = "Adversary"
$table[26,2] = 2
$table[26,3] = 16
$table[26,4] = 4
$table[27,0] = "Sumo (Remix) II"
$table[27,1] = "Adversary"
$table[27,2] = 2
$table[27,3] = 16
$table[27,4] = 4
$table[28,0] = "Sumo (Remix) III"
$table[28,1] = "Adversary"
$table[28,2] = 2
$table[28,3] = 16
$table[28,4] = 4
$table[29,0] = "Sumo (Remix) IV"
$table[29,1] = "Adversary"
$table[29,2] = 2
$table[29,3] = 16
$table[29,4] = 4
$table[30,0] = "Sumo (Remix) V"
$table[30,1] = "Adversary"
$table[30,2] = 2
$table[30,3] = 16
$table[30,4] = 4
$table[31,0] = "Sumo (Remix) VI"
$table[31,1] = "Adversary"
$table[31,2] = 2
$table[31,3] = 16
$table[31,4] = 4
$table[32,0] = "Sumo (Remix) VII"
$table[32,1] = "Adversary"
$table[32,2] = 2
$table[32,3] = 16
$table[32,4] = 4

$ws.Range("A1").Resize(33, 5).Value2 = $table

# --- Match the saved selection from the diff (cell B22 was selected on save).
$ws.Range("B22").Select() | Out-Null
